$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 15 & 16: Status column (G) moves from "In progress" to "Done" ---
$ws.Range("G15").Value = "Done"
$ws.Range("G12").Copy()
$ws.Range("G15").PasteSpecial(-4122)   # xlPasteFormats - reuse the existing "Done" (green) cell style

$ws.Range("G16").Value = "Done"
$ws.Range("G12").Copy()
$ws.Range("G16").PasteSpecial(-4122)   # xlPasteFormats - reuse the existing "Done" (green) cell style

# --- Row 27 / 28: assignees swapped (احمد <-> عبود) as part of the same housekeeping pass ---
$ws.Range("I27").Value = "عبود "
$ws.Range("I28").Value = "احمد"

# --- Rows 31-35: add assignee names in column I ---
$ws.Range("I31").Value = "احمد"
$ws.Range("I32").Value = "عبود"
$ws.Range("I33").Value = "عرين"
$ws.Range("I34").Value = "لمى"
$ws.Range("I35").Value = "احمد"

# --- Row 36: fill in Sprint, Status, Date and Assignee ---
$ws.Range("F36").Value = 5

$ws.Range("G36").Value = "In progress "
$ws.Range("G23").Copy()
$ws.Range("G36").PasteSpecial(-4122)   # xlPasteFormats - reuse the existing "In progress" (orange) cell style

$ws.Range("H36").Value = "20/04/2024-20/5/2024"
$ws.Range("H36").HorizontalAlignment = -4108   # xlCenter, matches style used by sibling columns

$ws.Range("I36").Value = "امجد وعمر"

$excel.CutCopyMode = $false

# --- View state: scroll / selection as left by the editor ---
$ws.Range("D38").Select()
$excel.ActiveWindow.Zoom = 100
